# Refresh the cryptocurrency ranking list (Price / Volume(1h) columns, and
# for rows that dropped off / entered the top-50, the Coin name + Link too).
# Mirrors the "Updated cryptos list ... with GitHub Actions" data refresh:
#   - row 4/7/10/11/... : Volume(1h) percentage only changes
#   - row 2/3/5/6/...    : Price + Volume(1h) change
#   - rows 35-51         : "USDe" dropped out of the list, every following
#                           coin shifted up one rank, and "SuiNetwork" is
#                           the new entry at the bottom (row 51) -> Coin,
#                           Link, Price and Volume(1h) all change there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="66.518.47"; E="  +3.93%  " },
    @{ Row=3; D="3.486.44"; E="  +2.59%  " },
    @{ Row=4; E="  +0.05%  " },
    @{ Row=5; D="590.37"; E="  +3.10%  " },
    @{ Row=6; D="167.83"; E="  +3.29%  " },
    @{ Row=7; E="  -0.03%  " },
    @{ Row=8; D="3.485.55"; E="  +2.62%  " },
    @{ Row=9; D="0.592"; E="  +7.61%  " },
    @{ Row=10; E="  +0.42%  " },
    @{ Row=11; E="  +5.75%  " },
    @{ Row=12; D="0.433"; E="  +3.09%  " },
    @{ Row=13; D="4.091.65"; E="  +2.67%  " },
    @{ Row=14; E="  -0.40%  " },
    @{ Row=15; D="28.03"; E="  +4.55%  " },
    @{ Row=16; D="66.554.85"; E="  +3.96%  " },
    @{ Row=17; E="  +2.60%  " },
    @{ Row=18; D="3.504.72"; E="  +3.03%  " },
    @{ Row=19; D="6.27"; E="  +2.62%  " },
    @{ Row=20; D="13.92"; E="  +3.94%  " },
    @{ Row=21; D="388.38"; E="  +4.03%  " },
    @{ Row=22; D="7.90"; E="  +1.64%  " },
    @{ Row=23; D="72.83"; E="  +3.79%  " },
    @{ Row=24; E="  -0.35%  " },
    @{ Row=25; E="  +3.75%  " },
    @{ Row=26; E="  +5.74%  " },
    @{ Row=27; D="10.27"; E="  +9.09%  " },
    @{ Row=28; E="  +1.34%  " },
    @{ Row=29; E="  +0.10%  " },
    @{ Row=30; D="6.29"; E="  +3.70%  " },
    @{ Row=31; D="1.44"; E="  +4.06%  " },
    @{ Row=32; E="  +2.56%  " },
    @{ Row=33; E="  +3.50%  " },
    @{ Row=34; D="7.34"; E="  +4.58%  " },
    @{ Row=35; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="1.59"; E="  +8.38%  " },
    @{ Row=36; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="162.76"; E="  +2.63%  " },
    @{ Row=37; B="Mantle"; C="https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D="0.888"; E="  +4.33%  " },
    @{ Row=38; B="Stacks"; C="https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D="1.91"; E="  +5.18%  " },
    @{ Row=39; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="6.75"; E="  +5.28%  " },
    @{ Row=40; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.0742"; E="  +2.76%  " },
    @{ Row=41; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="4.62"; E="  +6.49%  " },
    @{ Row=42; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="26.26"; E="  +2.27%  " },
    @{ Row=43; B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="2.779.11"; E="  +1.58%  " },
    @{ Row=44; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="43.06"; E="  +1.03%  " },
    @{ Row=45; B="InjectiveProtocol"; C="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D="26.47"; E="  +2.78%  " },
    @{ Row=46; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.0309"; E="  +2.05%  " },
    @{ Row=47; B="dogwifhat"; C="https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D="2.50"; E="  +4.18%  " },
    @{ Row=48; B="Bittensor"; C="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D="345.31"; E="  +5.33%  " },
    @{ Row=49; B="ONDO"; C="https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"; D="1.08"; E="  +4.33%  " },
    @{ Row=50; B="Arweave"; C="https://coinranking.com/coin/7XWg41D1+arweave-ar"; D="33.64"; E="  +12.21%  " },
    @{ Row=51; B="SuiNetwork"; C="https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"; D="0.864"; E="  +6.57%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($col in 'B','C','D','E') {
        if ($u.ContainsKey($col)) {
            $colIndex = switch ($col) { 'B' {2}; 'C' {3}; 'D' {4}; 'E' {5} }
            $cell = $ws.Cells.Item($row, $colIndex)
            # Force text so values such as "7.90", "1.00" or "66.518.47"
            # round-trip verbatim instead of Excel coercing them to numbers
            # (which would drop the formatted trailing/leading zeros), then
            # drop the forced format again so no stray style index sticks.
            $cell.NumberFormat = '@'
            $cell.Value = $u[$col]
            $cell.Style = 'Normal'
        }
    }
}
